$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing row 24 edits ---
# Text for C24 gets extended
$ws.Range("C24").Value = "presentatie voorbereiden, powerpoint uitbreiden"
# D24 hours increase from 1 to 2
$ws.Range("D24").Value = 2

# --- New rows 25-29: set the values first (so aggregate formulas over the
#     new range pick them up correctly), THEN copy the formatting used by
#     the existing log rows (row 24) onto the new rows. ---

# Row 25 - woensdag, laatste voorbereiding presentatie
$ws.Range("A25").Value = 42354
$ws.Range("B25").Value = "woensdag"
$ws.Range("C25").Value = "laatste voorbereiding presentatie"
$ws.Range("D25").Value = 1

# Row 26 - woensdag, tussenpresentatie
$ws.Range("A26").Value = 42354
$ws.Range("B26").Value = "woensdag"
$ws.Range("C26").Value = "tussenpresentatie"
$ws.Range("D26").Value = 1

# Row 27 - vrijdag, projectdocument bijwerken
$ws.Range("A27").Value = 42356
$ws.Range("B27").Value = "vrijdag"
$ws.Range("C27").Value = "projectdocument bijwerken"
$ws.Range("D27").Value = 1

# Row 28 - maandag, gebruikshandleiding maken
$ws.Range("A28").Value = 42359
$ws.Range("B28").Value = "maandag"
$ws.Range("C28").Value = "gebruikshandleiding maken"
$ws.Range("D28").Value = 1

# Row 29 - maandag, Globale opzet eindschermActivity en selectionActivity
$ws.Range("A29").Value = 42359
$ws.Range("B29").Value = "maandag"
$ws.Range("C29").Value = "Globale opzet eindschermActivity en selectionActivity"
$ws.Range("D29").Value = 1

$ws.Range("A24:D24").Copy()
$ws.Range("A25:D29").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

for ($r = 25; $r -le 29; $r++) {
    $ws.Rows.Item($r).RowHeight = $ws.Rows.Item(24).RowHeight
}

# --- Update the view/selection to reflect the new scroll position & selection ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("C5").Select()
